$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells L1:N1, styled to match the existing header row (A1:K1) ---
$ws.Range("L1").Value = "StrikePrice"
$ws.Range("M1").Value = "ExpiryDate"
$ws.Range("N1").Value = "OptionType"
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)   # xlPasteFormats - reuse header style

# --- Row 2 (HDFCBANK) - Tag text changes; L2:N2 are new trailing columns (left blank) ---
$ws.Range("K2").Value = "Stock Order"

# --- Row 3 -> SENSEX PUT option order ---
$ws.Range("A3").Value = "BSXOPT"
$ws.Range("B3").Value = "BFO"
$ws.Range("C3").Value = "BUY"
$ws.Range("D3").Value = 60
$ws.Range("E3").Value = "LIMIT"
$ws.Range("F3").Value = 138
$ws.Range("G3").Value = "INTRADAY"
$ws.Range("H3").Value = 238
$ws.Range("I3").Value = 90
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = "SENSEX PUT 85000"
$ws.Range("L3").Value = 85000
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "2025-12-18"
$ws.Range("M3").Style = "Normal"
$ws.Range("N3").Value = "PE"

# --- Row 4 -> SENSEX PUT option order ---
$ws.Range("A4").Value = "BSXOPT"
$ws.Range("B4").Value = "BFO"
$ws.Range("C4").Value = "BUY"
$ws.Range("D4").Value = 60
$ws.Range("E4").Value = "LIMIT"
$ws.Range("F4").Value = 140
$ws.Range("G4").Value = "INTRADAY"
$ws.Range("H4").Value = 240
$ws.Range("I4").Value = 92
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = "SENSEX PUT 84900"
$ws.Range("L4").Value = 84900
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = "2025-12-18"
$ws.Range("M4").Style = "Normal"
$ws.Range("N4").Value = "PE"

# --- Row 5 -> NIFTY call option order ---
$ws.Range("A5").Value = "NIFTY-Dec2025-24000-CE"
$ws.Range("B5").Value = "NSE_FNO"
$ws.Range("C5").Value = "BUY"
$ws.Range("D5").Value = 50
$ws.Range("E5").Value = "LIMIT"
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = "INTRADAY"
$ws.Range("H5").Value = 150
$ws.Range("I5").Value = 80
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = "NIFTY Call"
